# Refresh the cryptocurrency price table ("Updated cryptos list ... with GitHub Actions"):
#  - rows 33/34 (Filecoin <-> InternetComputer(DFINITY)) swapped position/rank
#  - Price (D) and Volume(1h) (E) columns refreshed with the latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (safe from Excel numeric auto-conversion)
$plainUpdates = @{
    'D2' = '30.334.12'
    'E2' = '  +2.15%  '
    'D3' = '2.091.14'
    'E3' = '  -0.21%  '
    'E4' = '  -0.86%  '
    'E5' = '  -0.60%  '
    'E6' = '  -0.77%  '
    'E7' = '  +1.52%  '
    'E8' = '  +0.66%  '
    'E9' = '  +3.81%  '
    'E10' = '  +0.98%  '
    'E11' = '  -0.27%  '
    'E12' = '  -0.03%  '
    'E13' = '  +3.67%  '
    'E14' = '  +2.27%  '
    'D15' = '2.036.21'
    'E15' = '  -2.48%  '
    'E16' = '  +1.86%  '
    'E17' = '  +0.68%  '
    'E18' = '  -0.75%  '
    'E19' = '  +1.31%  '
    'E20' = '  +0.12%  '
    'E21' = '  +2.26%  '
    'E22' = '  -0.68%  '
    'D23' = '30.317.56'
    'E23' = '  +1.92%  '
    'E24' = '  -0.30%  '
    'E25' = '  -0.77%  '
    'E26' = '  -0.45%  '
    'E27' = '  +0.41%  '
    'E28' = '  -0.53%  '
    'E29' = '  -0.06%  '
    'E30' = '  +0.28%  '
    'E31' = '  +1.43%  '
    'E32' = '  -0.58%  '
    'B33' = 'InternetComputer(DFINITY)'
    'C33' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E33' = '  +9.25%  '
    'B34' = 'Filecoin'
    'C34' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E34' = '  +1.07%  '
    'E35' = '  -1.86%  '
    'E36' = '  -0.66%  '
    'E37' = '  +2.49%  '
    'E38' = '  +2.09%  '
    'E39' = '  +1.80%  '
    'E40' = '  +3.43%  '
    'E41' = '  +0.44%  '
    'E42' = '  -0.89%  '
    'E43' = '  +2.73%  '
    'E44' = '  +0.78%  '
    'E45' = '  +1.45%  '
    'E47' = '  +18.21%  '
    'E48' = '  +0.26%  '
    'E49' = '  -1.17%  '
    'E50' = '  +7.90%  '
    'E51' = '  -0.37%  '
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Price (column D) updates whose text looks like a plain number; force text format
# so Excel keeps significant trailing/leading zeros instead of coercing to a double.
$priceTextUpdates = @{
    'D5' = '342.95'
    'D6' = '1.001'
    'D7' = '0.5235'
    'D8' = '0.4423'
    'D9' = '54.54'
    'D10' = '0.09329'
    'D11' = '1.168'
    'D12' = '24.85'
    'D13' = '8.590'
    'D14' = '6.896'
    'D16' = '101.21'
    'D17' = '0.00001160'
    'D18' = '1.002'
    'D19' = '21.12'
    'D20' = '0.06666'
    'D21' = '6.333'
    'D22' = '1.001'
    'D24' = '12.53'
    'D25' = '2.301'
    'D26' = '21.81'
    'D27' = '162.76'
    'D28' = '2.515'
    'D29' = '133.04'
    'D30' = '1.137'
    'D31' = '1.674'
    'D32' = '0.1045'
    'D33' = '6.778'
    'D34' = '6.244'
    'D35' = '3.862'
    'D36' = '10.21'
    'D37' = '0.02633'
    'D38' = '0.06841'
    'D39' = '0.6984'
    'D40' = '1.344'
    'D41' = '12.51'
    'D42' = '0.2208'
    'D43' = '0.6824'
    'D44' = '14.37'
    'D45' = '2.347'
    'D47' = '1.372'
    'D48' = '3.635'
    'D49' = '0.00000000348'
    'D50' = '1.210'
    'D51' = '1.216'
}
foreach ($addr in $priceTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $priceTextUpdates[$addr]
    $cell.Style = 'Normal'
}

Write-Output "Updated $($plainUpdates.Count + $priceTextUpdates.Count) cells in cryptos sheet"
